# Scheduled market-data refresh for Tonberry_Profits workbook.
# For each affected (sheet,row) pair, columns H-N -- currentAveragePrice,
# currentAveragePriceNQ, currentAveragePriceHQ, LevePriceNQ, LevePriceHQ,
# LeveProfitNQ, LeveProfitHQ -- are refreshed with the latest market values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 408.5
$ws.Range("I55").Value = 353.8
$ws.Range("J55").Value = 499.66666
$ws.Range("K55").Value = 353.8
$ws.Range("L55").Value = 499.66666
$ws.Range("M55").Value = -139.8
$ws.Range("N55").Value = -927.66666
$ws.Range("H88").Value = 1588.7778
$ws.Range("I88").Value = 2500
$ws.Range("J88").Value = 1474.875
$ws.Range("K88").Value = 2500
$ws.Range("L88").Value = 1474.875
$ws.Range("M88").Value = -2094
$ws.Range("N88").Value = -2286.875
$ws.Range("H91").Value = 1588.7778
$ws.Range("I91").Value = 2500
$ws.Range("J91").Value = 1474.875
$ws.Range("K91").Value = 2500
$ws.Range("L91").Value = 1474.875
$ws.Range("M91").Value = -1096
$ws.Range("N91").Value = -4282.875
$ws.Range("H100").Value = 1768.091
$ws.Range("I100").Value = 1744.9
$ws.Range("K100").Value = 1744.9
$ws.Range("M100").Value = -1203.9
$ws.Range("H113").Value = 22545.4
$ws.Range("J113").Value = 2000
$ws.Range("L113").Value = 2000
$ws.Range("N113").Value = -8508
$ws.Range("H132").Value = 1255.325
$ws.Range("I132").Value = 1051.6857
$ws.Range("K132").Value = 3155.0571
$ws.Range("M132").Value = -625.0571
$ws.Range("H138").Value = 2814.5908
$ws.Range("J138").Value = 3564.8
$ws.Range("L138").Value = 10694.4
$ws.Range("N138").Value = -20974.4
$ws.Range("H141").Value = 877562.9399999999
$ws.Range("I141").Value = 1121695
$ws.Range("K141").Value = 3365085
$ws.Range("M141").Value = -3359905

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5815103
$ws.Range("I2").Value = 11628706
$ws.Range("K2").Value = 11628706
$ws.Range("M2").Value = -11628593
$ws.Range("H45").Value = 1479.3158
$ws.Range("I45").Value = 907.1429000000001
$ws.Range("K45").Value = 907.1429000000001
$ws.Range("M45").Value = -530.1429000000001
$ws.Range("H88").Value = 3718.0908
$ws.Range("I88").Value = 2066.5
$ws.Range("J88").Value = 5700
$ws.Range("K88").Value = 2066.5
$ws.Range("L88").Value = 5700
$ws.Range("M88").Value = -1660.5
$ws.Range("N88").Value = -6512
$ws.Range("H91").Value = 3718.0908
$ws.Range("I91").Value = 2066.5
$ws.Range("J91").Value = 5700
$ws.Range("K91").Value = 2066.5
$ws.Range("L91").Value = 5700
$ws.Range("M91").Value = -662.5
$ws.Range("N91").Value = -8508
$ws.Range("H110").Value = 3928.25
$ws.Range("I110").Value = 1900
$ws.Range("K110").Value = 1900
$ws.Range("M110").Value = 145
$ws.Range("H116").Value = 5815103
$ws.Range("I116").Value = 11628706
$ws.Range("K116").Value = 11628706
$ws.Range("M116").Value = -11626412
$ws.Range("H122").Value = 1997.5
$ws.Range("I122").Value = 1997.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5992.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3542.5
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5815103
$ws.Range("I3").Value = 11628706
$ws.Range("K3").Value = 11628706
$ws.Range("M3").Value = -11628592
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("N49").ClearContents()
$ws.Range("H82").Value = 24662.666
$ws.Range("I82").Value = 11994
$ws.Range("K82").Value = 11994
$ws.Range("M82").Value = -11611
$ws.Range("H85").Value = 24662.666
$ws.Range("I85").Value = 11994
$ws.Range("K85").Value = 11994
$ws.Range("M85").Value = -10668
$ws.Range("H86").Value = 127575.25
$ws.Range("I86").Value = 1966.5
$ws.Range("J86").Value = 202940.5
$ws.Range("K86").Value = 1966.5
$ws.Range("L86").Value = 202940.5
$ws.Range("M86").Value = -843.5
$ws.Range("N86").Value = -205186.5
$ws.Range("H89").Value = 127575.25
$ws.Range("I89").Value = 1966.5
$ws.Range("J89").Value = 202940.5
$ws.Range("K89").Value = 9832.5
$ws.Range("L89").Value = 1014702.5
$ws.Range("M89").Value = -4216.5
$ws.Range("N89").Value = -1025934.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2067.9412
$ws.Range("I31").Value = 1364.1428
$ws.Range("K31").Value = 1364.1428
$ws.Range("M31").Value = -1069.1428
$ws.Range("H34").Value = 2067.9412
$ws.Range("I34").Value = 1364.1428
$ws.Range("K34").Value = 1364.1428
$ws.Range("M34").Value = -1162.1428
$ws.Range("H58").Value = 2417910
$ws.Range("I58").Value = 4349626
$ws.Range("K58").Value = 4349626
$ws.Range("M58").Value = -4349423
$ws.Range("H59").Value = 25166.666
$ws.Range("J59").Value = 25166.666
$ws.Range("L59").Value = 25166.666
$ws.Range("N59").Value = -27456.666
$ws.Range("H132").Value = 3400.8462
$ws.Range("I132").Value = 1466.5
$ws.Range("K132").Value = 4399.5
$ws.Range("M132").Value = -1869.5
$ws.Range("H134").Value = 1231.1482
$ws.Range("I134").Value = 1271.7084
$ws.Range("J134").Value = 906.6667
$ws.Range("K134").Value = 3815.1252
$ws.Range("L134").Value = 2720.0001
$ws.Range("M134").Value = -1280.1252
$ws.Range("N134").Value = -7790.0001
$ws.Range("H136").Value = 2417910
$ws.Range("I136").Value = 4349626
$ws.Range("K136").Value = 13048878
$ws.Range("M136").Value = -13046328

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1454
$ws.Range("I3").Value = 692.5
$ws.Range("K3").Value = 2077.5
$ws.Range("M3").Value = -1965.5
$ws.Range("H23").Value = 262.66666
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 262.66666
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 787.9999799999999
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -1257.99998
$ws.Range("H128").Value = 399999
$ws.Range("I128").Value = 399999
$ws.Range("K128").Value = 1199997
$ws.Range("M128").Value = -1195017
$ws.Range("H131").Value = 792.39
$ws.Range("J131").Value = 801.433
$ws.Range("L131").Value = 2404.299
$ws.Range("N131").Value = -12484.299
$ws.Range("H137").Value = 3489.4736
$ws.Range("I137").Value = 1827.5
$ws.Range("J137").Value = 3932.6667
$ws.Range("K137").Value = 5482.5
$ws.Range("L137").Value = 11798.0001
$ws.Range("M137").Value = -382.5
$ws.Range("N137").Value = -21998.0001
$ws.Range("H140").Value = 1476.079
$ws.Range("I140").Value = 793.9524
$ws.Range("J140").Value = 2318.7058
$ws.Range("K140").Value = 2381.8572
$ws.Range("L140").Value = 6956.117400000001
$ws.Range("M140").Value = 2798.1428
$ws.Range("N140").Value = -17316.1174

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1615.7
$ws.Range("I122").Value = 1209
$ws.Range("J122").Value = 2371
$ws.Range("K122").Value = 3627
$ws.Range("L122").Value = 7113
$ws.Range("M122").Value = -1177
$ws.Range("N122").Value = -12013

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4799.125
$ws.Range("I7").Value = 3380.8
$ws.Range("K7").Value = 3380.8
$ws.Range("M7").Value = -3268.8
$ws.Range("H22").Value = 1927.3572
$ws.Range("I22").Value = 1856.4286
$ws.Range("J22").Value = 1998.2858
$ws.Range("K22").Value = 1856.4286
$ws.Range("L22").Value = 1998.2858
$ws.Range("M22").Value = -1561.4286
$ws.Range("N22").Value = -2588.2858
$ws.Range("H27").Value = 1927.3572
$ws.Range("I27").Value = 1856.4286
$ws.Range("J27").Value = 1998.2858
$ws.Range("K27").Value = 1856.4286
$ws.Range("L27").Value = 1998.2858
$ws.Range("M27").Value = -1749.4286
$ws.Range("N27").Value = -2212.2858
$ws.Range("H46").Value = 2448.2
$ws.Range("I46").Value = 1424.75
$ws.Range("K46").Value = 1424.75
$ws.Range("M46").Value = -1236.75
$ws.Range("H55").Value = 319.54544
$ws.Range("I55").Value = 252.04347
$ws.Range("K55").Value = 252.04347
$ws.Range("M55").Value = -79.04347000000001
$ws.Range("H126").Value = 4799.125
$ws.Range("I126").Value = 3380.8
$ws.Range("K126").Value = 10142.4
$ws.Range("M126").Value = -7672.400000000001
$ws.Range("H136").Value = 2489.1155
$ws.Range("I136").Value = 1485.1052
$ws.Range("J136").Value = 5214.2856
$ws.Range("K136").Value = 4455.3156
$ws.Range("L136").Value = 15642.8568
$ws.Range("M136").Value = -1905.3156
$ws.Range("N136").Value = -20742.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1269.2222
$ws.Range("I132").Value = 990.4375
$ws.Range("K132").Value = 2971.3125
$ws.Range("M132").Value = -441.3125
$ws.Range("H135").Value = 82271.336
$ws.Range("J135").Value = 82271.336
$ws.Range("L135").Value = 82271.336
$ws.Range("N135").Value = -92411.336
